$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = -14.3013
$ws.Range("D4").Value = -7.631299999999999

# Row 5
$ws.Range("D5").Value = -8.074499999999999

# Row 7
$ws.Range("C7").Value = -11.7982

# Row 8
$ws.Range("D8").Value = -8.392999999999995

# Row 16
$ws.Range("C16").Value = -12.06870000000001
$ws.Range("D16").Value = -8.567100000000009
